$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Remove the "Meta description" paragraph that follows the H1 title.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2) Insert a new bold paragraph ("Play Fairy's Treasure Online for Free - Detailed
#    Slot Review") right before the final paragraph (the one currently holding the
#    italic image-prompt text), then give that final paragraph's italic run the new
#    meta-description copy. Both steps are done with InsertXML so straight
#    apostrophes are preserved verbatim (Find/Replace would smart-quote them).
$count = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs.Item($count - 1)
$beforeLast.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$newHeadingXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fairy''s Treasure Online for Free - Detailed Slot Review</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newHeadingXml)

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$metaDescXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover Fairy''s Treasure online slot machine for free with this detailed review. Find out about the game''s features and bonuses and play for free.</w:t></w:r></w:p>'
$lastPara.Range.InsertXML($metaDescXml)
